{"js": "// The document contains a single table of two-digit-number / one-digit-number\n// division problems laid out 5-per-row, with data rows interleaved with\n// blank spacer rows. We need to overwrite the text of 25 specific data\n// cells (in reading order: top-to-bottom, left-to-right) with new problem\n// strings. Because a couple of the \"before\" strings repeat (e.g.\n// \"22\u00f76=3, 4\" occurs twice, with two different replacements), we cannot do\n// a blind global find/replace by text \u2014 we must walk the table cells in\n// document order and apply the Nth replacement to the Nth non-blank cell.\n\nconst replacements = [\n  \"81\u00f79=9, 0\",\n  \"33\u00f73=11, 0\",\n  \"91\u00f75=18, 1\",\n  \"41\u00f72=20, 1\",\n  \"42\u00f72=21, 0\",\n  \"30\u00f75=6, 0\",\n  \"52\u00f77=7, 3\",\n  \"97\u00f78=12, 1\",\n  \"92\u00f76=15, 2\",\n  \"79\u00f77=11, 2\",\n  \"93\u00f75=18, 3\",\n  \"10\u00f75=2, 0\",\n  \"10\u00f79=1, 1\",\n  \"81\u00f79=9, 0\",\n  \"66\u00f74=16, 2\",\n  \"64\u00f78=8, 0\",\n  \"12\u00f77=1, 5\",\n  \"62\u00f79=6, 8\",\n  \"42\u00f72=21, 0\",\n  \"19\u00f78=2, 3\",\n  \"19\u00f73=6, 1\",\n  \"12\u00f79=1, 3\",\n  \"70\u00f72=35, 0\",\n  \"60\u00f77=8, 4\",\n  \"83\u00f75=16, 3\",\n];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\ntable.load(\"rowCount,values\");\nawait context.sync();\n\nlet next = 0;\nfor (let r = 0; r < table.values.length && next < replacements.length; r++) {\n  const row = table.values[r];\n  const rowHasData = row.some((cellText) => cellText !== \"\");\n  if (!rowHasData) continue;\n\n  for (let c = 0; c < row.length && next < replacements.length; c++) {\n    const cell = table.getCell(r, c);\n    cell.value = replacements[next];\n    next++;\n  }\n}\n\nawait context.sync();\n", "ps1": "# The document contains a single table of two-digit-number / one-digit-number\n# division problems laid out 5-per-row, with data rows interleaved with\n# blank spacer rows. We need to overwrite the text of 25 specific data\n# cells (in reading order: top-to-bottom, left-to-right) with new problem\n# strings. Because a couple of the \"before\" strings repeat (e.g.\n# \"22\u00f76=3, 4\" occurs twice, with two different replacements), we cannot do\n# a blind global find/replace by text \u2014 we must walk the table cells in\n# document order and apply the Nth replacement to the Nth non-blank cell.\n\n$replacements = @(\n  \"81\u00f79=9, 0\",\n  \"33\u00f73=11, 0\",\n  \"91\u00f75=18, 1\",\n  \"41\u00f72=20, 1\",\n  \"42\u00f72=21, 0\",\n  \"30\u00f75=6, 0\",\n  \"52\u00f77=7, 3\",\n  \"97\u00f78=12, 1\",\n  \"92\u00f76=15, 2\",\n  \"79\u00f77=11, 2\",\n  \"93\u00f75=18, 3\",\n  \"10\u00f75=2, 0\",\n  \"10\u00f79=1, 1\",\n  \"81\u00f79=9, 0\",\n  \"66\u00f74=16, 2\",\n  \"64\u00f78=8, 0\",\n  \"12\u00f77=1, 5\",\n  \"62\u00f79=6, 8\",\n  \"42\u00f72=21, 0\",\n  \"19\u00f78=2, 3\",\n  \"19\u00f73=6, 1\",\n  \"12\u00f79=1, 3\",\n  \"70\u00f72=35, 0\",\n  \"60\u00f77=8, 4\",\n  \"83\u00f75=16, 3\"\n)\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n$next = 0\nfor ($r = 1; $r -le $t.Rows.Count; $r++) {\n  if ($next -ge $replacements.Count) { break }\n\n  $rowHasData = $false\n  for ($c = 1; $c -le $t.Columns.Count; $c++) {\n    if ($t.Cell($r, $c).Range.Text.TrimEnd(\"`r\", \"`a\") -ne \"\") {\n      $rowHasData = $true\n      break\n    }\n  }\n  if (-not $rowHasData) { continue }\n\n  for ($c = 1; $c -le $t.Columns.Count; $c++) {\n    if ($next -ge $replacements.Count) { break }\n    $t.Cell($r, $c).Range.Text = $replacements[$next]\n    $next++\n  }\n}\n"}
